# Test multipart, fix error when mixing normal and multipart
#
# Insert three new columns (Multi Step Index, Block Index, Background)
# right after the existing "Pre or Post" column (old column D), pushing
# the previously-existing Question Stem..Machine TEKS columns from D:M
# to G:P. Then populate the three new columns for the header row and
# add two new "multipart" rows (4 and 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 blank columns at D (old D..M shift to G..P)
$ws.Columns("D:F").Insert()

# 2) Column widths for the new columns (values chosen so that, after this
#    engine's internal pixel-rounding on save, the stored widths match
#    14.5 / 10.33203125 / 16.83203125 as closely as the engine allows)
$ws.Columns("D").ColumnWidth = 13.666666666666666
$ws.Columns("E").ColumnWidth = 9.5
$ws.Columns("F").ColumnWidth = 16

# 3) Header row (row 1) for the new columns
$ws.Range("D1").Value = "Multi Step Index"
$ws.Range("E1").Value = "Block Index"
$ws.Range("F1").Value = "Background"

# 4) New data rows for the multipart question
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Oh look, multipart!"
$ws.Range("G4").Value = "Second part?"
$ws.Range("N4").Value = "Q3"
$ws.Range("O4").Value = "T1.3"
$ws.Range("P4").Value = "e6beb10a-f5cd-4b18-bf4a-e2f7174779bd"

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "Oh look, multipart!"
$ws.Range("G5").Value = "First part?"
$ws.Range("N5").Value = "Q3"
$ws.Range("O5").Value = "T1.3"
$ws.Range("P5").Value = "e6beb10a-f5cd-4b18-bf4a-e2f7174779bd"

# 5) Update the sheet view/selection to match the author's final state
$ws.Range("E6").Select()
